$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price/volume figures (GitHub Actions data pull).
# D/E columns are stored as text in the sheet; a leading apostrophe forces
# Excel to keep numeric-looking strings ("238.61", "1.0000", ...) as text
# instead of auto-converting them to numbers (matches the original inlineStr cells).

$ws.Range("D2").Value = "'25.730.15"
$ws.Range("E2").Value = "'  -2.74%  "

$ws.Range("D3").Value = "'1.741.57"
$ws.Range("E3").Value = "'  -5.02%  "

$ws.Range("D5").Value = "'238.61"
$ws.Range("E5").Value = "'  -8.24%  "

$ws.Range("E6").Value = "'  -0.04%  "

$ws.Range("D7").Value = "'0.5049"
$ws.Range("E7").Value = "'  -6.21%  "

$ws.Range("D8").Value = "'41.81"
$ws.Range("E8").Value = "'  -6.72%  "

$ws.Range("E9").Value = "'  -12.27%  "

$ws.Range("D10").Value = "'0.06146"
$ws.Range("E10").Value = "'  -10.36%  "

$ws.Range("D11").Value = "'1.744.75"
$ws.Range("E11").Value = "'  -5.49%  "

$ws.Range("D12").Value = "'0.06928"
$ws.Range("E12").Value = "'  -3.91%  "

$ws.Range("D13").Value = "'15.36"
$ws.Range("E13").Value = "'  -12.06%  "

$ws.Range("E14").Value = "'  -9.30%  "

$ws.Range("D15").Value = "'0.5950"
$ws.Range("E15").Value = "'  -19.03%  "

$ws.Range("D16").Value = "'76.57"
$ws.Range("E16").Value = "'  -13.91%  "

$ws.Range("D17").Value = "'1.0000"
$ws.Range("E17").Value = "'  -0.13%  "

$ws.Range("E18").Value = "'  -0.05%  "

$ws.Range("D19").Value = "'25.741.16"
$ws.Range("E19").Value = "'  -2.83%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000006809"
$ws.Range("E20").Value = "'  -13.46%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'11.61"
$ws.Range("E21").Value = "'  -15.86%  "

$ws.Range("D22").Value = "'1.967.07"
$ws.Range("E22").Value = "'  -5.47%  "

$ws.Range("D23").Value = "'4.042"
$ws.Range("E23").Value = "'  -11.42%  "

$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'5.174"
$ws.Range("E24").Value = "'  -13.05%  "

$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "'8.115"
$ws.Range("E25").Value = "'  -12.01%  "

$ws.Range("D26").Value = "'138.04"

$ws.Range("D27").Value = "'1.517"
$ws.Range("E27").Value = "'  -9.88%  "

$ws.Range("E28").Value = "'  -17.33%  "

$ws.Range("D29").Value = "'14.97"
$ws.Range("E29").Value = "'  -11.36%  "

$ws.Range("D30").Value = "'103.12"
$ws.Range("E30").Value = "'  -6.54%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'3.756"
$ws.Range("E31").Value = "'  -10.68%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08090"
$ws.Range("E32").Value = "'  -8.08%  "

$ws.Range("D33").Value = "'3.461"
$ws.Range("E33").Value = "'  -13.33%  "

$ws.Range("D34").Value = "'0.04493"
$ws.Range("E34").Value = "'  -6.13%  "

$ws.Range("D35").Value = "'0.9996"
$ws.Range("E35").Value = "'  -0.05%  "

$ws.Range("D36").Value = "'2.651"
$ws.Range("E36").Value = "'  -9.77%  "

$ws.Range("D37").Value = "'0.9777"
$ws.Range("E37").Value = "'  -13.20%  "

$ws.Range("D38").Value = "'0.6080"
$ws.Range("E38").Value = "'  -16.29%  "

$ws.Range("D39").Value = "'2.650"
$ws.Range("E39").Value = "'  -14.17%  "

$ws.Range("D40").Value = "'0.01548"
$ws.Range("E40").Value = "'  -9.10%  "

$ws.Range("D41").Value = "'1.000"
$ws.Range("E41").Value = "'  -0.01%  "

$ws.Range("D42").Value = "'1.896"
$ws.Range("E42").Value = "'  -16.54%  "

$ws.Range("D43").Value = "'103.09"
$ws.Range("E43").Value = "'  -4.21%  "

$ws.Range("D44").Value = "'0.3797"
$ws.Range("E44").Value = "'  -19.31%  "

$ws.Range("D45").Value = "'5.097"
$ws.Range("E45").Value = "'  -13.15%  "

$ws.Range("D46").Value = "'0.7320"
$ws.Range("E46").Value = "'  -19.06%  "

$ws.Range("D47").Value = "'0.05342"
$ws.Range("E47").Value = "'  -7.68%  "

$ws.Range("D48").Value = "'0.1112"
$ws.Range("E48").Value = "'  -9.55%  "

$ws.Range("D49").Value = "'30.11"

$ws.Range("D50").Value = "'5.885"
$ws.Range("E50").Value = "'  -19.62%  "

$ws.Range("D51").Value = "'52.48"
$ws.Range("E51").Value = "'  -12.32%  "
